$d = $word.ActiveDocument

# Locate the existing run containing "今天天气不错" and place the insertion
# point right after it, then append the new sentence so it becomes a new
# run with the same character formatting (as in the diff, a separate <w:r>).
$target = $d.Content
$target.Find.Execute("今天天气不错", $false, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null

$target.Collapse(0)
$insertStart = $target.End
$target.InsertAfter("，心情也很好。")
$newRange = $d.Range($insertStart, $target.End)
$newRange.Font.NameFarEast = "+Body"
$newRange.Font.Size = 24
$newRange.Font.SizeBi = 24

